$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 1212, shifting the existing data (rows 1212-1302)
# down to 1213-1303, and populate the new row with the new price entry.
$ws.Rows("1212:1212").Insert()

$ws.Range("A1212").Value = 8
$ws.Range("B1212").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1212").Value = "Coquimbo"
$ws.Range("D1212").Value = 45265
$ws.Range("E1212").Value = 4
$ws.Range("F1212").Value = 100112023
$ws.Range("G1212").Value = "Brócoli"
$ws.Range("H1212").Value = "Sin especificar"
$ws.Range("I1212").Value = "Primera"
$ws.Range("J1212").Value = 2100
$ws.Range("K1212").Value = 800
$ws.Range("L1212").Value = 900
$ws.Range("M1212").Value = 850
$ws.Range("N1212").Value = "$/unidad"
$ws.Range("O1212").Value = "Provincia del Elquí"
$ws.Range("P1212").Value = 850
$ws.Range("Q1212").Value = 1
$ws.Range("R1212").Value = "Hortaliza"
